# Added grades for weeks 6 and 7
$wb = $excel.ActiveWorkbook

$noteText = "Just a tip, usually you would validate the input before performing any action. Also, though I haven't tested it, I think there is a minor bug in your gridOn code. I think that if you don't supply a gridOn parameter, this will call gridOff since you don't have else with the exist condition. Obviously no points off, just thought you might want some pointers. Kudos on using the short-circuit AND."

# ---- Week 6 ----
$ws6 = $wb.Worksheets.Item("Week 6")

# Narrow the "Grade" column and the "Notes" column to match the other graded weeks.
$ws6.Columns.Item(2).ColumnWidth = 8.33
$ws6.Columns.Item(4).ColumnWidth = 5.83

for ($r = 2; $r -le 6; $r++) {
    $ws6.Cells.Item($r, 2).Value = 1
    $ws6.Cells.Item($r, 4).ClearFormats()
}

# ---- Week 7 ----
$ws7 = $wb.Worksheets.Item("Week 7")

$ws7.Columns.Item(2).ColumnWidth = 8.33
$ws7.Columns.Item(4).ColumnWidth = 357.33

for ($r = 2; $r -le 14; $r++) {
    $ws7.Cells.Item($r, 2).Value = 1
    if ($r -eq 11) {
        $ws7.Cells.Item($r, 4).Value = $noteText
    } else {
        $ws7.Cells.Item($r, 4).ClearFormats()
    }
}
